# Adds a run of text "Test" to the last non-empty paragraph of the
# document body -- the paragraph immediately preceding the trailing
# empty paragraph / sectPr. That paragraph currently has no runs,
# only a <w:pPr><w:rPr><w:lang w:eastAsia="en-AU"/></w:rPr></w:pPr>.

$d = $word.ActiveDocument

# The very last paragraph in the body is an empty trailing paragraph
# (the one right before the section properties). The target paragraph
# is the one right before it.
$count  = $d.Paragraphs.Count
$target = $d.Paragraphs.Item($count - 1)

$r = $target.Range
# Exclude the paragraph mark from the range so the new text is
# inserted as a run inside the paragraph (not after it).
$r.End = $r.End - 1

$r.InsertAfter("Test")
# Match the paragraph's existing run formatting (w:lang w:eastAsia="en-AU").
$r.LanguageIDFarEast = "en-AU"
